# Weekly price update: insert a new daily record for row 151 ("Poroto verde",
# Terminal La Palmera de La Serena) and push the existing historical rows
# down by one, extending the data range from A1:R232 to A1:R233.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 151 - this shifts rows 151..232 down to
# 152..233 and grows the sheet dimension automatically.
$ws.Rows(151).Insert()

# The row that used to be 151 is now 152; duplicate its (unchanged)
# categorical fields into the new row 151 so every column keeps a sane
# value, then overwrite the columns that actually hold the new record's data.
$ws.Range("A152:R152").Copy()
$ws.Range("A151").PasteSpecial()

$ws.Range("D151").Value = 44719
$ws.Range("J151").Value = 520
$ws.Range("K151").Value = 19000
$ws.Range("L151").Value = 20000
$ws.Range("M151").Value = 19500
$ws.Range("O151").Value = "Perú"
$ws.Range("P151").Value = 780
